$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Add the new row 10 data (order matters for shared-string append order)
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = 200
$ws.Range("C10").Value = "Wait Stats"
$ws.Range("D10").Value = "(One per wait type)"
$ws.Range("E10").Value = "http://BrentOzar.com/waits/(waittype)"

# Update the title cell (A1) with the new date (appended last to shared strings)
$ws.Range("A1").Value = "sp_AskBrent Check ID List - v1 July 11, 2013"

# Update the active selection to A2
$ws.Range("A2").Select()
